# Natmi following Dr Hou advice:
# Recompute the Ncam1->Robo3 LR-pair table (rows 2-4 updated, two new
# rows added for "Neutro" and "sCs" sending clusters) with revised
# specificity statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 5,20
$data[0,0] = "ECs"
$data[0,1] = "Ncam1"
$data[0,2] = "Robo3"
$data[0,3] = "ECs"
$data[0,4] = 2
$data[0,5] = 0.6666666666666666
$data[0,6] = 0.6816986666666667
$data[0,7] = 2.045096
$data[0,8] = 0.01110711092851045
$data[0,9] = 0.01110711092851045
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 1.935833
$data[0,13] = 5.807499
$data[0,14] = 1
$data[0,15] = 1
$data[0,16] = 1.319654774989333
$data[0,17] = 11.876892974904
$data[0,18] = 0.01110711092851045
$data[0,19] = 0.01110711092851045
$data[1,0] = "FAPs"
$data[1,1] = "Ncam1"
$data[1,2] = "Robo3"
$data[1,3] = "ECs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 8.469728666666667
$data[1,7] = 25.409186
$data[1,8] = 0.1379997063732729
$data[1,9] = 0.1379997063732729
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 1.935833
$data[1,13] = 5.807499
$data[1,14] = 1
$data[1,15] = 1
$data[1,16] = 16.39598025397933
$data[1,17] = 147.563822285814
$data[1,18] = 0.1379997063732729
$data[1,19] = 0.1379997063732729
$data[2,0] = "M1"
$data[2,1] = "Ncam1"
$data[2,2] = "Robo3"
$data[2,3] = "ECs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 0.3748256666666667
$data[2,7] = 1.124477
$data[2,8] = 0.006107141559886993
$data[2,9] = 0.006107141559886992
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 1.935833
$data[2,13] = 5.807499
$data[2,14] = 1
$data[2,15] = 1
$data[2,16] = 0.7255998947803333
$data[2,17] = 6.530399053022999
$data[2,18] = 0.006107141559886993
$data[2,19] = 0.006107141559886992
$data[3,0] = "Neutro"
$data[3,1] = "Ncam1"
$data[3,2] = "Robo3"
$data[3,3] = "ECs"
$data[3,4] = 2
$data[3,5] = 0.6666666666666666
$data[3,6] = 0.1449283333333333
$data[3,7] = 0.434785
$data[3,8] = 0.002361358696634494
$data[3,9] = 0.002361358696634494
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 1.935833
$data[3,13] = 5.807499
$data[3,14] = 1
$data[3,15] = 1
$data[3,16] = 0.2805570503016667
$data[3,17] = 2.525013452715
$data[3,18] = 0.002361358696634494
$data[3,19] = 0.002361358696634494
$data[4,0] = "sCs"
$data[4,1] = "Ncam1"
$data[4,2] = "Robo3"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 51.70379466666666
$data[4,7] = 155.111384
$data[4,8] = 0.8424246824416953
$data[4,9] = 0.8424246824416952
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 1.935833
$data[4,13] = 5.807499
$data[4,14] = 1
$data[4,15] = 1
$data[4,16] = 100.0899119409573
$data[4,17] = 900.8092074686159
$data[4,18] = 0.8424246824416953
$data[4,19] = 0.8424246824416952
$ws.Range("A2:T6").Value = $data
